# Referral Follow-up form change: remove the "role" question from the
# "person-create" survey (select_one roles / role / Role) and blank out
# the now-unused "roles" choice list entries in the choices sheet.

$wb = $excel.ActiveWorkbook

# --- survey sheet: drop the "role" field (row 25) -------------------------
# Deleting the row shifts every subsequent row up by one, which also moves
# along each cell's existing style, exactly matching the target layout.
$wsSurvey = $wb.Worksheets.Item("survey")
$wsSurvey.Rows("25:25").Delete()

# --- choices sheet: clear the "roles" choice list (rows 20-25) ------------
# The list_name/name/label values are removed, but the row/cell styling is
# left in place (matches the diff: cells keep their `s` but lose `t`/`v`).
$wsChoices = $wb.Worksheets.Item("choices")
$wsChoices.Range("A20:C25").ClearContents()
